$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17. This shifts the old row 17 (and everything
# below it, down to the old last row 77) down by one, so old row 17 becomes
# new row 18, ..., old row 77 becomes new row 78. The new row 17 is blank and
# ready to receive the new weekly price entry.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new data point (most of the descriptive
# columns mirror the entry that used to be in row 17 / is now in row 18).
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 44622
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112021
$ws.Range("G17").Value = "Ají"
$ws.Range("H17").Value = "Chilena(o)"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = 21000
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = 21571
$ws.Range("N17").Value = "$/saco 25 kilos"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 863
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
